# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The worker/period data table (rows 16-36, columns C-G) is reorganized:
# previously grouped by worker (7 consecutive rows per worker, periods
# descending 2308->2302), now grouped by period (3 consecutive rows per
# period block, periods ascending 2302->2308).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Worker master data: Doc Number (C), Name (D), Valor Mora F-column value,
# Salario Basico G-column value
$workers = @(
    @("1026272073", "GABRIEL ANTONIO GARRIDO ROMERO", 60000, 1500000),
    @("1047493694", "JOSE DAVID FERNANDEZ RHENALS", 48000, 1200000),
    @("1047505194", "LUIS MARIO CASTRILLO ARRIETA", 100000, 2500000)
)

$periods = @("2302", "2303", "2304", "2305", "2306", "2307", "2308")

$row = 16
foreach ($period in $periods) {
    foreach ($worker in $workers) {
        $doc = $worker[0]
        $name = $worker[1]
        $colFValue = $worker[2]
        $colGValue = $worker[3]

        $ws.Cells.Item($row, 3).Value = $doc        # C: N Doc Trabajador
        $ws.Cells.Item($row, 4).Value = $name        # D: Nombre Trabajador
        $ws.Cells.Item($row, 5).Value = $period      # E: Periodo Mora
        $ws.Cells.Item($row, 6).Value = $colFValue   # F: Valor Mora
        $ws.Cells.Item($row, 7).Value = $colGValue   # G: Salario Basico

        $row = $row + 1
    }
}
